# Update cryptocurrency price/volume data pulled on 24-1-2023 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'302.64"
$ws.Range("E2").Value = "'-1.14%"
$ws.Range("D3").Value = "'35.49"
$ws.Range("E3").Value = "'-2.32%"
$ws.Range("D4").Value = "'5.083"
$ws.Range("E4").Value = "'0.33%"
$ws.Range("D5").Value = "'0.08063"
$ws.Range("E5").Value = "'1.74%"
$ws.Range("D6").Value = "'1.942"
$ws.Range("E6").Value = "'-11.82%"
$ws.Range("D7").Value = "'7.844"
$ws.Range("E7").Value = "'-2.15%"
$ws.Range("D8").Value = "'2.898"
$ws.Range("E8").Value = "'10.09%"
$ws.Range("D9").Value = "'0.9227"
$ws.Range("E9").Value = "'-0.86%"
$ws.Range("D10").Value = "'0.1080"
$ws.Range("E10").Value = "'9.80%"
$ws.Range("D11").Value = "'0.1896"
$ws.Range("E11").Value = "'1.00%"
$ws.Range("D12").Value = "'0.09520"
$ws.Range("E12").Value = "'4.41%"
$ws.Range("D13").Value = "'0.03678"
$ws.Range("E13").Value = "'0.14%"
$ws.Range("D14").Value = "'0.09908"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("D15").Value = "'0.001454"
$ws.Range("E15").Value = "'1.24%"
$ws.Range("D16").Value = "'0.005796"
$ws.Range("E16").Value = "'3.13%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.41%"
$ws.Range("E18").Value = "'-0.82%"
$ws.Range("D19").Value = "'0.3416"
$ws.Range("E19").Value = "'1.31%"
$ws.Range("E20").Value = "'-2.49%"
$ws.Range("D21").Value = "'5.151"
$ws.Range("E21").Value = "'0.98%"
$ws.Range("E22").Value = "'0.33%"
$ws.Range("D23").Value = "'0.04542"
$ws.Range("E23").Value = "'-0.38%"
$ws.Range("D24").Value = "'0.001231"
$ws.Range("E24").Value = "'-0.62%"
$ws.Range("D25").Value = "'0.004716"
$ws.Range("E25").Value = "'-1.40%"
$ws.Range("D26").Value = "'0.0001262"
$ws.Range("E26").Value = "'-2.83%"
$ws.Range("D27").Value = "'0.0004465"
$ws.Range("E27").Value = "'-5.78%"
$ws.Range("D39").Value = "'0.01936"
$ws.Range("E39").Value = "'0.11%"
$ws.Range("D40").Value = "'0.04776"
$ws.Range("E40").Value = "'-3.18%"
$ws.Range("D41").Value = "'0.007587"
$ws.Range("E41").Value = "'-2.87%"
$ws.Range("D42").Value = "'0.009649"
$ws.Range("E42").Value = "'23.57%"
$ws.Range("D44").Value = "'0.002133"
$ws.Range("E44").Value = "'1.03%"
$ws.Range("E45").Value = "'0.10%"
$ws.Range("D46").Value = "'0.00006491"
$ws.Range("E46").Value = "'4.00%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("D48").Value = "'64.42"
$ws.Range("E48").Value = "'24.12%"
$ws.Range("D49").Value = "'0.001305"
$ws.Range("E49").Value = "'-27.53%"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.30%"
